# Replace Product, IT, and Finance templates with correct industry-specific
# content - this workbook is the Finance Change Management Plan template,
# which is being re-themed from an "AI/ML" initiative to a "Banking"
# implementation project.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Change Management Overview" ---
$ws1 = $wb.Worksheets.Item("Change Management Overview")

# Title / header block
$ws1.Range("A2").Value = "Banking Implementation Project"

# Project information
$ws1.Range("B6").Value = "Enterprise Banking Implementation"

# Change management objectives
$ws1.Range("A15").Value = "1. Achieve 95% user adoption of new Banking systems within 6 months of go-live"
$ws1.Range("A17").Value = "3. Build organizational capability and confidence in Banking technologies"
$ws1.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for Banking transformation"

# Touch these previously-blank rows (no cell content) so their row metadata
# gets written out, matching the original authoring tool's row bookkeeping.
$ws1.Rows.Item(13).OutlineLevel = 0
$ws1.Rows.Item(21).OutlineLevel = 0

# --- Sheet 2: "Change Impact Assessment" ---
$ws2 = $wb.Worksheets.Item("Change Impact Assessment")

$ws2.Range("G4").Value = "Banking automation"

$ws2.Rows.Item(2).OutlineLevel = 0

# --- Sheet 3: "Change Activities" ---
$ws3 = $wb.Worksheets.Item("Change Activities")

$ws3.Rows.Item(2).OutlineLevel = 0
